# Swap the presentation's theme (currently the "Integral" / Red Violet
# design used by the slide master) for the plain "Office Theme" colour
# scheme that was otherwise sitting unused on the notes master.
#
# Only the colour scheme actually differs between the two theme parts in
# this deck (font scheme and format scheme are already identical), so the
# edit is expressed as 12 RGB writes on the slide master's theme colour
# scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) -- the
# documented PowerPoint COM surface for editing a theme's colours
# (Theme.ThemeColorScheme.Item(i).RGB / .Colors(i).RGB).

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$cs = $theme.ThemeColorScheme

# VBA/COM RGB() packs colour bytes as 0xBBGGRR, so build the value from
# the target hex (RRGGBB) accordingly.
function BGR([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Index order of ThemeColorScheme.Item(i): 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1..accent6, 11 hlink, 12 folHlink.
$officeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $cs.Item($i).RGB = BGR $officeColors[$i - 1]
}
